$wb = $excel.ActiveWorkbook

$values = @{
    "C11" = 0.39348498093710305
    "D11" = -0.3000000000000007
    "E11" = 0.6135817094968132
    "F11" = -0.17400000000000015
    "G11" = 1.5829618029997903
    "H11" = 16.12947350163202
    "I11" = 1.65096649881347
}

foreach ($sheetName in @("Test 1", "Test 2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }
}
